# Workaround for Word Online shortcoming: reorder <w:ilvl/> before
# <w:numId/> inside <w:numPr> for every list paragraph in the document
# (matches upstream commit "Change order of ilvl and numId in
# document.xml (#5647)").
#
# Re-assigning ListFormat.ListLevelNumber to its own current value forces
# Word to re-emit the paragraph's <w:numPr> block, which is enough to
# normalize the element order without altering the paragraph's actual
# list membership or level.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $lf = $p.Range.ListFormat
    if ($lf.ListType -ne 0) {
        $lf.ListLevelNumber = $lf.ListLevelNumber
    }
}

Write-Output "Reordered ilvl/numId for $($d.Paragraphs.Count) paragraphs checked"
